$d = $word.ActiveDocument

# The thesis paragraph originally reads (single run):
#   "...that can repair the body better than itself. This is advanced form of
#   treatment is attributed to Stem Cells...cannot. "
# The author finished editing there, which is where Word drops its "last
# edit" bookmark. Word keeps only one "_GoBack" bookmark in a document, so
# the one currently sitting near " VII.)" must be removed and a new one
# added around "This is advanced " in the thesis paragraph.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("This is advanced ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$d.Bookmarks.Add("_GoBack", $rng)
